$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 59
$ws1.Range("F4").Value = 1476
$ws1.Range("F5").Value = 504
$ws1.Range("F6").Value = 1059
$ws1.Range("F7").Value = 10966
$ws1.Range("F8").Value = 10966
$ws1.Range("F11").Value = 311
$ws1.Range("F13").Value = 745
$ws1.Range("F14").Value = 12192
$ws1.Range("F15").Value = 12694
$ws1.Range("F22").Value = 15

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 59
$ws4.Range("F5").Value = 1476
$ws4.Range("F6").Value = 504
$ws4.Range("F7").Value = 1059
$ws4.Range("F8").Value = 10966
$ws4.Range("F9").Value = 10966
$ws4.Range("F12").Value = 311
$ws4.Range("F14").Value = 745
$ws4.Range("F15").Value = 12192
$ws4.Range("F16").Value = 12694
$ws4.Range("F23").Value = 15
